# ExcelImportSeeder-DatabaseSeeder update table relations
#
# - Drop the two unused, empty worksheets (Sheet2, Sheet3).
# - Rename the third column of the "Table_3__2" query table from
#   "field_id" to "sub_field_name", and refresh column C's data so it
#   references the name of the parent field instead of a numeric id.
# - Leave the active selection on the new header cell (C1), matching
#   the last user interaction captured in the workbook.

$wb = $excel.ActiveWorkbook

# --- Remove the empty, unused sheets -------------------------------------
$wb.Worksheets.Item("Sheet2").Delete()
$wb.Worksheets.Item("Sheet3").Delete()

# --- Update the "field_id" column on Sheet1 into "sub_field_name" --------
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C1").Value = "sub_field_name"

$ws.Range("C2").ClearContents()
$ws.Range("C3").ClearContents()
$ws.Range("C4").Value = "mobile"
$ws.Range("C5").Value = "cloud"
$ws.Range("C6").ClearContents()
$ws.Range("C7").Value = "network"
$ws.Range("C8").ClearContents()

# --- Match the workbook's final selection state ---------------------------
$ws.Range("C1").Select()
